# Updated symbol list with GitHub Actions - apply price/volume refresh
# and coin reordering per the upstream commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "306.80"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "36.19"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.90%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.061"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.12%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08064"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.08%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.146"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "11.42%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.840"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.02%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9264"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.36%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1419"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "10.69%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1921"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.65%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09076"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.64%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03449"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.94%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09910"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.34%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001401"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.34%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.006319"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-4.69%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.840"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "6.37%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.148"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.21%"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.398"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "11.79%"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3449"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.79%"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1336"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.02%"
$ws.Range("B21").Value = "MCDex"
$ws.Range("C21").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.799"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-7.07%"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2342"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-7.47%"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04363"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.05%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001230"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.33%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004300"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-9.00%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02007"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "0.26%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05149"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.25%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007505"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.64%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01013"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.28%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1361"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.41%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002150"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009959"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-6.94%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006279"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.70%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.02%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "64.85"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-0.16%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001250"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-21.89%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.02%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.02%"
